# Refresh the "cryptos" price/volume snapshot (Price = column D,
# Volume(1h) = column E) for rows 2-51 with the latest scraped values.
#
# Numeric-looking Price strings (e.g. "6.74") are written with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# inline-string / text cell type) instead of auto-converting them to
# numbers; the quote-prefix mark is then cleared via Style = "Normal" so
# the cell's style stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.318.44"
$ws.Range("E2").Value = "  +4.74%  "
$ws.Range("D3").Value = "3.257.00"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'577.89"
$ws.Range("E5").Value = "  +2.24%  "
$ws.Range("D6").Value = "'179.68"
$ws.Range("E6").Value = "  +5.81%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "3.255.67"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("E10").Value = "  +3.95%  "
$ws.Range("D11").Value = "'6.74"
$ws.Range("E11").Value = "  +1.44%  "
$ws.Range("D12").Value = "'0.415"
$ws.Range("E12").Value = "  +4.43%  "
$ws.Range("D13").Value = "3.824.49"
$ws.Range("E13").Value = "  +2.62%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "'28.23"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").Value = "67.281.73"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("E17").Value = "  +3.02%  "
$ws.Range("D18").Value = "3.257.09"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "'5.88"
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "'13.41"
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("D21").Value = "'376.49"
$ws.Range("E21").Value = "  +6.70%  "
$ws.Range("E22").Value = "  +6.28%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'71.40"
$ws.Range("E24").Value = "  +3.39%  "
$ws.Range("D25").Value = "'0.513"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").Value = "3.396.88"
$ws.Range("E26").Value = "  +2.49%  "
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'9.90"
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +3.94%  "
$ws.Range("D32").Value = "'5.64"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'22.59"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("E35").Value = "  +5.86%  "
$ws.Range("D36").Value = "'6.82"
$ws.Range("E36").Value = "  +2.70%  "
$ws.Range("D37").Value = "'163.76"
$ws.Range("E37").Value = "  +6.06%  "
$ws.Range("E38").Value = "  +4.25%  "
$ws.Range("D39").Value = "'0.857"
$ws.Range("E39").Value = "  +4.61%  "
$ws.Range("E40").Value = "  +9.83%  "
$ws.Range("D41").Value = "'27.01"
$ws.Range("E41").Value = "  +4.40%  "
$ws.Range("D42").Value = "'2.63"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").Value = "'6.65"
$ws.Range("E43").Value = "  +10.29%  "
$ws.Range("D44").Value = "2.773.73"
$ws.Range("E44").Value = "  +5.65%  "
$ws.Range("E45").Value = "  +5.47%  "
$ws.Range("D46").Value = "'25.84"
$ws.Range("E46").Value = "  +8.72%  "
$ws.Range("D47").Value = "'345.99"
$ws.Range("E47").Value = "  +7.04%  "
$ws.Range("D48").Value = "'40.49"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").Value = "'0.0675"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("D50").Value = "'0.0281"
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("E51").Value = "  +1.10%  "

# Reset style for cells written with a leading apostrophe (quote-prefix)
# so they do not pick up an extra quotePrefix style flag.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
